# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the latest generated output (gh-pages commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 17114
$ws1.Range("F4").Value  = 54
$ws1.Range("F6").Value  = 78
$ws1.Range("F8").Value  = 1054
$ws1.Range("F9").Value  = 406
$ws1.Range("F10").Value = 241
$ws1.Range("F12").Value = 11920
$ws1.Range("F14").Value = 59
$ws1.Range("F15").Value = 11634
$ws1.Range("F16").Value = 4722
$ws1.Range("F17").Value = 514
$ws1.Range("F18").Value = 62
$ws1.Range("F24").Value = 48

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 17114
$ws4.Range("F4").Value  = 54
$ws4.Range("F6").Value  = 78
$ws4.Range("F8").Value  = 1054
$ws4.Range("F9").Value  = 406
$ws4.Range("F10").Value = 241
$ws4.Range("F14").Value = 11920
$ws4.Range("F16").Value = 59
$ws4.Range("F17").Value = 11634
$ws4.Range("F18").Value = 4722
$ws4.Range("F19").Value = 514
$ws4.Range("F20").Value = 62
$ws4.Range("F26").Value = 48

$wb.Save()
